$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing rows (old 2013-2020 data, rows 15-22) that are no longer needed;
# this also shrinks the sheet dimension down to A1:G14.
$ws.Rows("15:22").Delete()

# Overwrite rows 2-14 (previously 2000-2012) with the refreshed series (2010-2022).
$ws.Range("A2").Value = "2010年"
$ws.Range("B2").Value = 1481.78395769851
$ws.Range("C2").Value = 2073.0625039957
$ws.Range("D2").Value = 2064.1833817007
$ws.Range("E2").Value = 409.7005925772
$ws.Range("F2").Value = 2783.9992197635
$ws.Range("G2").Value = 3199.3321450138
$ws.Range("A3").Value = "2011年"
$ws.Range("B3").Value = 1614.46016741067
$ws.Range("C3").Value = 2271.0572247187
$ws.Range("D3").Value = 2249.8919799959
$ws.Range("E3").Value = 426.8026803796
$ws.Range("F3").Value = 3048.148669316
$ws.Range("G3").Value = 3541.2983665006
$ws.Range("A4").Value = "2012年"
$ws.Range("B4").Value = 1729.64421193479
$ws.Range("C4").Value = 2449.6471794721
$ws.Range("D4").Value = 2443.9574475871
$ws.Range("E4").Value = 445.8660011663
$ws.Range("F4").Value = 3292.4303448605
$ws.Range("G4").Value = 3837.4798275499
$ws.Range("A5").Value = "2013年"
$ws.Range("B5").Value = 1851.59682678184
$ws.Range("C5").Value = 2639.8904562937
$ws.Range("D5").Value = 2618.4226828954
$ws.Range("E5").Value = 462.8410108002
$ws.Range("F5").Value = 3565.6746858027
$ws.Range("G5").Value = 4143.9875647812
$ws.Range("A6").Value = "2014年"
$ws.Range("B6").Value = 1976.59369057032
$ws.Range("C6").Value = 2835.92248236
$ws.Range("D6").Value = 2839.5228798828
$ws.Range("E6").Value = 481.5882389252
$ws.Range("F6").Value = 3863.1167407192
$ws.Range("G6").Value = 4440.5133060675
$ws.Range("A7").Value = "2015年"
$ws.Range("B7").Value = 2103.50556266842
$ws.Range("C7").Value = 3035.6091110882
$ws.Range("D7").Value = 3021.12420613865
$ws.Range("E7").Value = 500.3630597988
$ws.Range("F7").Value = 4202.1712796005
$ws.Range("G7").Value = 4703.6920396017
$ws.Range("A8").Value = "2016年"
$ws.Range("B8").Value = 2234.72677162222
$ws.Range("C8").Value = 3243.5107605796
$ws.Range("D8").Value = 3227.42778671986
$ws.Range("E8").Value = 516.7640280106
$ws.Range("F8").Value = 4542.3824267922
$ws.Range("G8").Value = 4987.3944281694
$ws.Range("A9").Value = "2017年"
$ws.Range("B9").Value = 2375.55619344151
$ws.Range("C9").Value = 3468.8439658698
$ws.Range("D9").Value = 3464.29860571328
$ws.Range("E9").Value = 537.2345549531
$ws.Range("F9").Value = 4918.0271569616
$ws.Range("G9").Value = 5280.0691721337
$ws.Range("A10").Value = "2018年"
$ws.Range("B10").Value = 2524.06885318935
$ws.Range("C10").Value = 3702.9830881698
$ws.Range("D10").Value = 3686.71881130701
$ws.Range("E10").Value = 555.9459912819
$ws.Range("F10").Value = 5310.6583480456
$ws.Range("G10").Value = 5585.9122593244
$ws.Range("A11").Value = "2019年"
$ws.Range("B11").Value = 2664.79368959606
$ws.Range("C11").Value = 3923.32912474348
$ws.Range("D11").Value = 3912.33681502289
$ws.Range("E11").Value = 573.051301730205
$ws.Range("F11").Value = 5692.02204626978
$ws.Range("G11").Value = 5858.20428911231
$ws.Range("A12").Value = "2020年"
$ws.Range("B12").Value = 2717.97119436732
$ws.Range("C12").Value = 4011.15827537261
$ws.Range("D12").Value = 3979.04070116392
$ws.Range("E12").Value = 590.994022566519
$ws.Range("F12").Value = 5802.85267032334
$ws.Range("G12").Value = 6002.56513920516
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 2945
$ws.Range("C13").Value = 4350
$ws.Range("D13").Value = 4319.7
$ws.Range("E13").Value = 632.8
$ws.Range("F13").Value = 6295.3
$ws.Range("G13").Value = 6523.2
$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 3033.4
$ws.Range("C14").Value = 4480.1
$ws.Range("D14").Value = 4432.1
$ws.Range("E14").Value = 659
$ws.Range("F14").Value = 6439.7
$ws.Range("G14").Value = 6768.3
